$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "45.019.76"
$ws.Range("D3").Value = "2.268.59"
Set-TextValue $ws.Range("D5") "301.16"
Set-TextValue $ws.Range("D6") "94.55"
Set-TextValue $ws.Range("D7") "0.566"
Set-TextValue $ws.Range("D9") "0.510"
Set-TextValue $ws.Range("D10") "34.12"
Set-TextValue $ws.Range("D12") "7.23"
$ws.Range("D14").Value = "2.612.57"
$ws.Range("D15").Value = "2.266.61"
Set-TextValue $ws.Range("D16") "13.62"
Set-TextValue $ws.Range("D17") "0.802"
$ws.Range("D18").Value = "44.914.85"
Set-TextValue $ws.Range("D19") "13.41"
$ws.Range("D20").Value = "0.0₃0917"
Set-TextValue $ws.Range("D21") "6.04"
Set-TextValue $ws.Range("D22") "65.63"
Set-TextValue $ws.Range("D23") "239.26"
Set-TextValue $ws.Range("D24") "2.88"
Set-TextValue $ws.Range("D25") "0.999"
Set-TextValue $ws.Range("D26") "1.90"
Set-TextValue $ws.Range("D27") "41.45"
Set-TextValue $ws.Range("D30") "19.64"
Set-TextValue $ws.Range("D31") "152.46"
Set-TextValue $ws.Range("D32") "5.53"
Set-TextValue $ws.Range("D33") "0.0790"
Set-TextValue $ws.Range("D34") "2.55"
Set-TextValue $ws.Range("D35") "2.92"
Set-TextValue $ws.Range("D38") "1.77"
Set-TextValue $ws.Range("D39") "3.91"
Set-TextValue $ws.Range("D40") "0.0309"
Set-TextValue $ws.Range("D41") "3.22"
Set-TextValue $ws.Range("D42") "13.70"
$ws.Range("D45").Value = "1.769.44"
Set-TextValue $ws.Range("D47") "76.67"
Set-TextValue $ws.Range("D48") "69.76"
Set-TextValue $ws.Range("D49") "95.71"
Set-TextValue $ws.Range("D50") "53.76"
Set-TextValue $ws.Range("D51") "7.89"

$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  -0.93%  "
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("E10").Value = "  -2.03%  "
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("E19").Value = "  +12.38%  "
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("E25").Value = "  -0.42%  "
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("E27").Value = "  +10.99%  "
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("E32").Value = "  -6.48%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  -4.03%  "
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("E37").Value = "  -5.06%  "
$ws.Range("E38").Value = "  -4.46%  "
$ws.Range("E39").Value = "  +3.71%  "
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("E42").Value = "  -9.64%  "
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("E44").Value = "  +12.12%  "
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("E46").Value = "  +2.45%  "
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  -1.40%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
